$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.085911989212036
$ws.Range("B1").Value = 2.158337116241455
$ws.Range("C1").Value = 2.235931396484375
$ws.Range("D1").Value = 2.976590156555176
$ws.Range("E1").Value = 2.99866247177124
